$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 19, pushing the existing rows 19-55
# down to become rows 21-57.
$ws.Rows("19:20").Insert()

# --- New row 19 ---
$ws.Cells.Item(19, 1).Value = 5
$ws.Cells.Item(19, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(19, 3).Value = "Maule"
$ws.Cells.Item(19, 4).Value = 44998
$ws.Cells.Item(19, 5).Value = 7
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100104
$ws.Cells.Item(19, 8).Value = "Frutos de pepita"
$ws.Cells.Item(19, 9).Value = 100104003
$ws.Cells.Item(19, 10).Value = "Membrillo"
$ws.Cells.Item(19, 11).Value = "Champion"
$ws.Cells.Item(19, 12).Value = "Especial"
$ws.Cells.Item(19, 13).Value = 200
$ws.Cells.Item(19, 14).Value = 12000
$ws.Cells.Item(19, 15).Value = 12000
$ws.Cells.Item(19, 16).Value = 12000
$ws.Cells.Item(19, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(19, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(19, 19).Value = 667
$ws.Cells.Item(19, 20).Value = 18

# --- New row 20 ---
$ws.Cells.Item(20, 1).Value = 5
$ws.Cells.Item(20, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(20, 3).Value = "Maule"
$ws.Cells.Item(20, 4).Value = 44998
$ws.Cells.Item(20, 5).Value = 7
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100104
$ws.Cells.Item(20, 8).Value = "Frutos de pepita"
$ws.Cells.Item(20, 9).Value = 100104003
$ws.Cells.Item(20, 10).Value = "Membrillo"
$ws.Cells.Item(20, 11).Value = "Champion"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 250
$ws.Cells.Item(20, 14).Value = 10000
$ws.Cells.Item(20, 15).Value = 10000
$ws.Cells.Item(20, 16).Value = 10000
$ws.Cells.Item(20, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(20, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(20, 19).Value = 556
$ws.Cells.Item(20, 20).Value = 18
